# Add newly-added ("some more paid courses") course links to the main
# "Sheet1" (Udemy Courses) worksheet, in column C only, starting at row 20
# and stepping by 2 (matching the existing sparse layout), then leave the
# selection/active sheet positioned the way the author left it when they
# saved: Sheet1 active (instead of the Kotlin sheet) with C42 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$urls = @(
    "https://www.udemy.com/course/digital-electronics-logic-design/",
    "https://www.udemy.com/course/mongodb-the-complete-developers-guide/",
    "https://www.udemy.com/course/burp-suite-mastery-bug-hunters-perspective/",
    "https://www.udemy.com/course/practical-ethical-hacking/",
    "https://www.udemy.com/course/linux-mastery/",
    "https://www.udemy.com/course/build-10-c-beginner-projects-from-scratch/",
    "https://www.udemy.com/course/linux-privilege-escalation-for-beginners/",
    "https://www.udemy.com/course/learn-python-and-ethical-hacking-from-scratch/",
    "https://www.udemy.com/course/the-complete-nmap-ethical-hacking-course-network-security/",
    "https://www.udemy.com/course/advanced-css-and-sass/",
    "https://www.udemy.com/course/responsive-web-design-tutorial-course-html5-css3-bootstrap/"
)

$row = 20
foreach ($url in $urls) {
    $ws.Cells.Item($row, 3).Value = $url
    $row += 2
}

# Activate Sheet1 and leave the selection/scroll position where the
# author left it before saving.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 5
$excel.ActiveWindow.ScrollColumn = 1
[void]$ws.Range("C42").Select()
